# Auto update Excel log
$wb = $excel.ActiveWorkbook

function Add-LogRow($ws, $row, $date, $timestamp, $hour, $location, $value, $status) {
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 6).NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value = $date
    $ws.Cells.Item($row, 2).Value = $timestamp
    $ws.Cells.Item($row, 3).Value = $hour
    $ws.Cells.Item($row, 4).Value = $location
    $ws.Cells.Item($row, 5).Value = $value
    $ws.Cells.Item($row, 6).Value = $status
}

# --- ALERTS sheet: add row 6 ---
$wsAlerts = $wb.Worksheets.Item("ALERTS")
Add-LogRow $wsAlerts 6 "2026-01-30" "15:44:42" "15:00" "Living Room" "CRITICAL EMERGENCY" "FALL_DETECTED"

# --- PIR sheet: add rows 78-89 ---
$wsPir = $wb.Worksheets.Item("PIR")
$pirTimestamps = @(
    "15:43:49", "15:43:53", "15:43:58", "15:44:03", "15:44:08", "15:44:13",
    "15:44:18", "15:44:23", "15:44:28", "15:44:33", "15:44:39", "15:44:43"
)
$row = 78
foreach ($ts in $pirTimestamps) {
    Add-LogRow $wsPir $row "2026-01-30" $ts "15:00" "Bathroom" "No Motion" "Inactive"
    $row++
}

# --- Humidity sheet: add rows 43-50 ---
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityTimestamps = @("15:43:49", "15:43:53", "15:44:04", "15:44:09", "15:44:14", "15:44:24", "15:44:29", "15:44:44")
$humidityValues = @("87.9%", "86.6%", "87.9%", "87.9%", "87.9%", "87.9%", "87.9%", "87.9%")
$row = 43
for ($i = 0; $i -lt $humidityTimestamps.Length; $i++) {
    Add-LogRow $wsHumidity $row "2026-01-30" $humidityTimestamps[$i] "15:00" "Bathroom" $humidityValues[$i] "Active"
    $row++
}

# --- mmWave sheet: add rows 6-9 ---
$wsMmwave = $wb.Worksheets.Item("mmWave")
$mmwaveTimestamps = @("15:43:48", "15:43:56", "15:44:06", "15:44:17")
$row = 6
foreach ($ts in $mmwaveTimestamps) {
    Add-LogRow $wsMmwave $row "2026-01-30" $ts "15:00" "Living Room" "PRESENCE_DETECTED" "Active"
    $row++
}
